$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new "No" column at the very left (shifts B..U -> C..V, merges/widths shift too) ---
$null = $ws.Columns("A:A").Insert()

# --- New header cell content ---
$ws.Range("A1").Value2 = "No"

# --- Borrow the existing header box-border/fill/font formatting for the two new header cells
#     (row1 "No" + row1 "Car Maker") from the neighbouring already-styled header cells, then
#     trim the shared inner edge so the two-row merge reads as one unbroken box. ---
$null = $ws.Range("C1").Copy()
$null = $ws.Range("A1:B1").PasteSpecial(-4122)
$null = $ws.Range("C2").Copy()
$null = $ws.Range("A2:B2").PasteSpecial(-4122)
$ws.Range("A1:B1").Borders.Item(9).LineStyle = -4142
$ws.Range("A2:B2").Borders.Item(8).LineStyle = -4142
$ws.Application.CutCopyMode = $false

# --- Merge the new "No" header cell across both header rows, like the other header cells ---
$null = $ws.Range("A1:A2").Merge()

# --- Column widths (re-tuned by the author after the insert) ---
$ws.Columns("A:A").ColumnWidth = 2.6666666666666665
$ws.Columns("C:C").ColumnWidth = 12
$ws.Columns("D:D").ColumnWidth = 7.666666666666667
$ws.Columns("E:F").ColumnWidth = 7
$ws.Columns("H:H").ColumnWidth = 11.333333333333334
$ws.Columns("I:I").ColumnWidth = 12.666666666666666
$ws.Columns("J:J").ColumnWidth = 10.666666666666666
$ws.Columns("L:M").ColumnWidth = 12
$ws.Columns("N:N").ColumnWidth = 6.333333333333333
$ws.Columns("O:O").ColumnWidth = 9.166666666666666
$ws.Columns("P:P").ColumnWidth = 14.833333333333334
$ws.Columns("R:R").ColumnWidth = 14.166666666666666
$ws.Columns("S:S").ColumnWidth = 14.833333333333334
$ws.Columns("U:U").ColumnWidth = 13.5

# --- Selection / view bookkeeping, matching what the author's Excel session left behind ---
$null = $ws.Columns("F:F").Select()
